$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(443).Insert()

$ws.Cells.Item(443, 1).Value = 10
$ws.Cells.Item(443, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(443, 3).Value = "La Araucanía"
$ws.Cells.Item(443, 4).Value = 45131
$ws.Cells.Item(443, 5).Value = 9
$ws.Cells.Item(443, 6).Value = 100112044
$ws.Cells.Item(443, 7).Value = "Perejil"
$ws.Cells.Item(443, 8).Value = "Sin especificar"
$ws.Cells.Item(443, 9).Value = "Primera"
$ws.Cells.Item(443, 10).Value = 65
$ws.Cells.Item(443, 11).Value = 4000
$ws.Cells.Item(443, 12).Value = 4000
$ws.Cells.Item(443, 13).Value = 4000
$ws.Cells.Item(443, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(443, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(443, 16).Value = 1333
$ws.Cells.Item(443, 17).Value = 3
$ws.Cells.Item(443, 18).Value = "Hortaliza"
